$wb = $excel.ActiveWorkbook

# ---- Sheet "Rushing" (sheet1) ----
$rushing = $wb.Worksheets.Item("Rushing")

# Week/game counter in column A shifts down by one for every data row (rows 2-11)
$rushing.Range("A2").Value = 0
$rushing.Range("A3").Value = 1
$rushing.Range("A4").Value = 2
$rushing.Range("A5").Value = 3
$rushing.Range("A6").Value = 4
$rushing.Range("A7").Value = 5
$rushing.Range("A8").Value = 6
$rushing.Range("A9").Value = 7
$rushing.Range("A10").Value = 8
$rushing.Range("A11").Value = 9

# S.Barkley (row 3) updated rushing attempt splits
$rushing.Range("C3").Value = 20
$rushing.Range("D3").Value = 7
$rushing.Range("E3").Value = 3

# D.Booker (row 4) updated rushing attempt splits
$rushing.Range("C4").Value = 50
$rushing.Range("D4").Value = 38
$rushing.Range("E4").Value = 9
$rushing.Range("F4").Value = 16

# ---- Sheet "Receiving" (sheet2) ----
$receiving = $wb.Worksheets.Item("Receiving")

# S.Barkley (row 2)
$receiving.Range("C2").Value = 20
$receiving.Range("D2").Value = 16

# D.Booker (row 3)
$receiving.Range("C3").Value = 27
$receiving.Range("D3").Value = 25

# K.Golladay (row 6)
$receiving.Range("C6").Value = 33
$receiving.Range("D6").Value = 26
$receiving.Range("E6").Value = 9
$receiving.Range("F6").Value = 6
$receiving.Range("G6").Value = 6

# D.Slayton (row 8)
$receiving.Range("C8").Value = 37
$receiving.Range("D8").Value = 25
$receiving.Range("E8").Value = 14
$receiving.Range("G8").Value = 5

# J.Ross (row 11)
$receiving.Range("C11").Value = 12

# P.Cooper (row 13)
$receiving.Range("C13").Value = 5
$receiving.Range("D13").Value = 3

# E.Engram (row 14)
$receiving.Range("C14").Value = 49
$receiving.Range("D14").Value = 33
$receiving.Range("E14").Value = 5

# K.Rudolph (row 15)
$receiving.Range("C15").Value = 25
$receiving.Range("D15").Value = 17
